$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.765.22'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.85%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.624.40'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5070'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2555'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06356'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07772'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.246'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.94%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.621.19'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.847.73'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5534'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.62'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₅7498'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.776.79'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.94%  '
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.06%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.393'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.761'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.974'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.002'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.36%  '
$ws.Range("E25").Value = '  -1.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1239'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.727'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.46'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.37%  '
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04861'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.310'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.83%  '
$ws.Range("E33").Value = '  -1.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.544'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.361'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8919'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.134.38'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5496'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.535'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01557'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.564'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7934'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.13'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.770.70'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.41%  '
$ws.Range("E46").Value = '  -6.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4417'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.63'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05126'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.619'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9980'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.81%  '
